# Auto-generated edit script applying Goblin_Profits market-data refresh
# Updates currentAveragePrice/LevePrice/LeveProfit columns (H-N) per row,
# matching the scheduled market-data refresh described in the commit.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 10000069
$ws.Range("I4").Value = 76.77778000000001
$ws.Range("K4").Value = 76.77778000000001
$ws.Range("M4").Value = 37.22221999999999

$ws.Range("H33").Value = 603.4783
$ws.Range("I33").Value = 163.58824
$ws.Range("K33").Value = 163.58824
$ws.Range("M33").Value = 65.41175999999999

$ws.Range("H98").Value = 7388.3335
$ws.Range("I98").Value = 9162.053
$ws.Range("K98").Value = 9162.053
$ws.Range("M98").Value = -7664.053

$ws.Range("H106").Value = 5854.1875
$ws.Range("I106").Value = 5474.846
$ws.Range("J106").Value = 7498
$ws.Range("K106").Value = 5474.846
$ws.Range("L106").Value = 7498
$ws.Range("M106").Value = -4843.846
$ws.Range("N106").Value = -8760

$ws.Range("H116").Value = 5900.778
$ws.Range("I116").Value = 5682.5
$ws.Range("J116").Value = 6337.3335
$ws.Range("K116").Value = 5682.5
$ws.Range("L116").Value = 6337.3335
$ws.Range("M116").Value = -2240.5
$ws.Range("N116").Value = -13221.3335

$ws.Range("H122").Value = 7388.3335
$ws.Range("I122").Value = 9162.053
$ws.Range("K122").Value = 27486.159
$ws.Range("M122").Value = -25036.159

$ws.Range("H137").Value = 1562.2222
$ws.Range("I137").Value = 1438.5454
$ws.Range("J137").Value = 1756.5714
$ws.Range("K137").Value = 4315.6362
$ws.Range("L137").Value = 5269.7142
$ws.Range("M137").Value = -1765.6362
$ws.Range("N137").Value = -10369.7142

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3242.3125
$ws.Range("I32").Value = 2555.7046
$ws.Range("K32").Value = 2555.7046
$ws.Range("M32").Value = -2268.7046

$ws.Range("H76").Value = 45999.75
$ws.Range("J76").Value = 45999.75
$ws.Range("L76").Value = 45999.75
$ws.Range("N76").Value = -46675.75

$ws.Range("H79").Value = 45999.75
$ws.Range("J79").Value = 45999.75
$ws.Range("L79").Value = 45999.75
$ws.Range("N79").Value = -48339.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1482.5834
$ws.Range("I22").Value = 500.5
$ws.Range("J22").Value = 1679
$ws.Range("K22").Value = 500.5
$ws.Range("L22").Value = 1679
$ws.Range("M22").Value = -150.5
$ws.Range("N22").Value = -2379

$ws.Range("H31").Value = 2652.24
$ws.Range("I31").Value = 1300.5
$ws.Range("K31").Value = 1300.5
$ws.Range("M31").Value = -1005.5

$ws.Range("H34").Value = 2652.24
$ws.Range("I34").Value = 1300.5
$ws.Range("K34").Value = 1300.5
$ws.Range("M34").Value = -1098.5

$ws.Range("H50").Value = 25066.2
$ws.Range("I50").Value = 25083
$ws.Range("J50").Value = 24999
$ws.Range("K50").Value = 25083
$ws.Range("L50").Value = 24999
$ws.Range("M50").Value = -24458
$ws.Range("N50").Value = -26249

$ws.Range("H58").Value = 1470.1111
$ws.Range("I58").Value = 1247.4286
$ws.Range("J58").Value = 2249.5
$ws.Range("K58").Value = 1247.4286
$ws.Range("L58").Value = 2249.5
$ws.Range("M58").Value = -1044.4286
$ws.Range("N58").Value = -2655.5

$ws.Range("H60").Value = 41664.5
$ws.Range("J60").Value = 41997.6
$ws.Range("L60").Value = 41997.6
$ws.Range("N60").Value = -43019.6

$ws.Range("H136").Value = 1470.1111
$ws.Range("I136").Value = 1247.4286
$ws.Range("J136").Value = 2249.5
$ws.Range("K136").Value = 3742.2858
$ws.Range("L136").Value = 6748.5
$ws.Range("M136").Value = -1192.2858
$ws.Range("N136").Value = -11848.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H50").Value = 14540.667
$ws.Range("I50").Value = 751.25
$ws.Range("K50").Value = 2253.75
$ws.Range("M50").Value = -1772.75

$ws.Range("H53").Value = 14540.667
$ws.Range("I53").Value = 751.25
$ws.Range("K53").Value = 2253.75
$ws.Range("M53").Value = -1772.75

$ws.Range("H63").Value = 200
$ws.Range("I63").Value = 200
$ws.Range("K63").Value = 600
$ws.Range("M63").Value = 149

$ws.Range("H66").Value = 200
$ws.Range("I66").Value = 200
$ws.Range("K66").Value = 1800
$ws.Range("M66").Value = 1944

$ws.Range("H80").Value = 19833
$ws.Range("J80").Value = 19666.666
$ws.Range("L80").Value = 58999.99800000001
$ws.Range("N80").Value = -60871.99800000001

$ws.Range("H82").Value = 27866.555
$ws.Range("I82").Value = 15819.8
$ws.Range("J82").Value = 42925
$ws.Range("K82").Value = 47459.39999999999
$ws.Range("L82").Value = 128775
$ws.Range("M82").Value = -47053.39999999999
$ws.Range("N82").Value = -129587

$ws.Range("H83").Value = 19833
$ws.Range("J83").Value = 19666.666
$ws.Range("L83").Value = 176999.994
$ws.Range("N83").Value = -186359.994

$ws.Range("H85").Value = 27866.555
$ws.Range("I85").Value = 15819.8
$ws.Range("J85").Value = 42925
$ws.Range("K85").Value = 47459.39999999999
$ws.Range("L85").Value = 128775
$ws.Range("M85").Value = -46055.39999999999
$ws.Range("N85").Value = -131583

$ws.Range("H131").Value = 1373.0834
$ws.Range("I131").Value = 661.0909
$ws.Range("K131").Value = 1983.2727
$ws.Range("M131").Value = 3056.7273

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 4034.3635
$ws.Range("I3").Value = 3486.4443
$ws.Range("J3").Value = 6500
$ws.Range("K3").Value = 3486.4443
$ws.Range("L3").Value = 6500
$ws.Range("M3").Value = -3370.4443
$ws.Range("N3").Value = -6732

$ws.Range("H10").Value = 722500
$ws.Range("J10").Value = 9583.333000000001
$ws.Range("L10").Value = 9583.333000000001
$ws.Range("N10").Value = -9921.333000000001

$ws.Range("H11").Value = 12886818
$ws.Range("I11").Value = 14172000
$ws.Range("K11").Value = 14172000
$ws.Range("M11").Value = -14171861

$ws.Range("H14").Value = 6333.5557
$ws.Range("I14").Value = 3667.6667
$ws.Range("J14").Value = 3667.6667
$ws.Range("K14").Value = 3667.6667
$ws.Range("L14").Value = 7666.5
$ws.Range("M14").Value = -3499.6667
$ws.Range("N14").Value = -8002.5

$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("N110").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4000
$ws.Range("I7").Value = 4000
$ws.Range("K7").Value = 4000
$ws.Range("M7").Value = -3888

$ws.Range("H42").Value = 999999.5
$ws.Range("I42").Value = 999999.5
$ws.Range("K42").Value = 999999.5
$ws.Range("M42").Value = -999436.5

$ws.Range("H46").Value = 3293.8
$ws.Range("I46").Value = 1225
$ws.Range("K46").Value = 1225
$ws.Range("M46").Value = -1037

$ws.Range("H49").Value = 999999.5
$ws.Range("I49").Value = 999999.5
$ws.Range("K49").Value = 999999.5
$ws.Range("M49").Value = -999852.5

$ws.Range("H68").Value = 6547.381
$ws.Range("I68").Value = 4916
$ws.Range("J68").Value = 7199.933
$ws.Range("K68").Value = 4916
$ws.Range("L68").Value = 7199.933
$ws.Range("M68").Value = -4167
$ws.Range("N68").Value = -8697.933000000001

$ws.Range("H71").Value = 6547.381
$ws.Range("I71").Value = 4916
$ws.Range("J71").Value = 7199.933
$ws.Range("K71").Value = 24580
$ws.Range("L71").Value = 35999.665
$ws.Range("M71").Value = -20836
$ws.Range("N71").Value = -43487.665

$ws.Range("H82").Value = 3042
$ws.Range("I82").Value = 1537
$ws.Range("J82").Value = 5299.5
$ws.Range("K82").Value = 1537
$ws.Range("L82").Value = 5299.5
$ws.Range("M82").Value = -1176
$ws.Range("N82").Value = -6021.5

$ws.Range("H85").Value = 3042
$ws.Range("I85").Value = 1537
$ws.Range("J85").Value = 5299.5
$ws.Range("K85").Value = 1537
$ws.Range("L85").Value = 5299.5
$ws.Range("M85").Value = -289
$ws.Range("N85").Value = -7795.5

$ws.Range("H126").Value = 4000
$ws.Range("I126").Value = 4000
$ws.Range("K126").Value = 12000
$ws.Range("M126").Value = -9530

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 20000
$ws.Range("I34").Value = 20000
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 20000
$ws.Range("L34").Value = 0
$ws.Range("N34").ClearContents()

$ws.Range("H123").Value = 33331.332
$ws.Range("J123").Value = 33331.332
$ws.Range("L123").Value = 33331.332
$ws.Range("N123").Value = -43131.332
